$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D7","D8","D9","D11","D13","D14","D15","D18","D21","D22","D24","D25","D26","D28","D29","D30","D32","D33","D35","D37","D41","D42","D43","D45","D46","D47","D48","D49","D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.291.17"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.677.77"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "217.44"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "0.5255"
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("D7").Value = "1.008"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "0.2686"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("D9").Value = "0.06463"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").Value = "0.07508"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "1.676.04"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "4.511"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "0.5766"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "0.000008467"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "26.324.45"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "4.916"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("D21").Value = "189.49"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "6.178"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "144.93"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "7.775"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "0.1262"
$ws.Range("E26").Value = "  +6.32%  "
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").Value = "0.06440"
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("D29").Value = "1.365"
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("D30").Value = "1.322"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "3.580"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").Value = "1.654"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").Value = "0.6181"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").Value = "2.740"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").Value = "1.114.49"
$ws.Range("E39").Value = "  +3.67%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "0.8711"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").Value = "1.015"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "100.44"
$ws.Range("D44").Value = "1.827.96"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").Value = "56.89"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.009"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "8.143"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000103"
$ws.Range("E48").Value = "  -11.14%  "
$ws.Range("D49").Value = "0.05265"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").Value = "6.045"
$ws.Range("E51").Value = "  +1.72%  "
